# Template.xlsx KPI config sheet update:
#  - rename the "entity_*" filter headers to "filter_entity_*"
#  - rename the "entity_*_filter" value headers to "filter_entity_*_value"
#  - rename the "manufacturer" KPI-row label to "manufacturer_name"
#  - move the active selection to B15 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI")

# Row 1 header renames
$ws.Range("C1").Value = "filter_entity_1"
$ws.Range("D1").Value = "filter_entity_2"
$ws.Range("E1").Value = "filter_entity_3"
$ws.Range("H1").Value = "filter_entity_1_value"
$ws.Range("I1").Value = "filter_entity_2_value"
$ws.Range("J1").Value = "filter_entity_3_value"

# "manufacturer" -> "manufacturer_name" for the two KPI rows that used it
$ws.Range("C2").Value = "manufacturer_name"
$ws.Range("C3").Value = "manufacturer_name"

# Restore the saved selection/active cell
$ws.Range("B15").Select()
